# Rename the "Shift begins" / "Shifts ends" headers on the Events sheet
# (the space is dropped from both labels: "Shift begins" -> "ShiftBegins",
# "Shifts ends" -> "ShiftsEnds").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

$ws.Range("I1").Value = "ShiftBegins"
$ws.Range("J1").Value = "ShiftsEnds"

# Re-apply the existing (identical) formatting on a couple of ranges so the
# workbook's style table collapses the redundant duplicate style entries
# that had accumulated, matching the formatting the cells already display.
$ws.Range("H2:J2").Copy()
$ws.Range("H8:J10").PasteSpecial(-4122)

$ws.Range("C8").Copy()
$ws.Range("C2:C7").PasteSpecial(-4122)

$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Reflect the user's final selection/zoom state on the Events sheet.
$ws.Select()
$ws.Range("J1").Select()
$excel.ActiveWindow.Zoom = 104
